$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 3425.8823
$ws.Range("I19").Value = 1819.1111
$ws.Range("K19").Value = 1819.1111
$ws.Range("M19").Value = -1644.1111
# Row 32
$ws.Range("H32").Value = 7442
$ws.Range("J32").Value = 7425.4287
$ws.Range("L32").Value = 7425.4287
$ws.Range("N32").Value = -8077.4287
# Row 62
$ws.Range("H62").Value = 6369.5835
$ws.Range("I62").Value = 5722
$ws.Range("K62").Value = 5722
$ws.Range("M62").Value = -5098
# Row 65
$ws.Range("H65").Value = 6369.5835
$ws.Range("I65").Value = 5722
$ws.Range("K65").Value = 28610
$ws.Range("M65").Value = -25490
# Row 92
$ws.Range("H92").Value = 1251.4615
$ws.Range("I92").Value = 1029
$ws.Range("J92").Value = 1607.4
$ws.Range("K92").Value = 1029
$ws.Range("L92").Value = 1607.4
$ws.Range("M92").Value = 219
$ws.Range("N92").Value = -4103.4

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1464.0667
$ws.Range("J2").Value = 1505.5
$ws.Range("L2").Value = 1505.5
$ws.Range("N2").Value = -1731.5
# Row 24
$ws.Range("H24").Value = 34999.5
$ws.Range("J24").Value = 34999.5
$ws.Range("L24").Value = 34999.5
$ws.Range("N24").Value = -35747.5
# Row 32
$ws.Range("H32").Value = 11499888
$ws.Range("I32").Value = 14494913
$ws.Range("J32").Value = 18958.166
$ws.Range("K32").Value = 14494913
$ws.Range("L32").Value = 18958.166
$ws.Range("M32").Value = -14494626
$ws.Range("N32").Value = -19532.166
# Row 45
$ws.Range("H45").Value = 1653.1428
$ws.Range("I45").Value = 1319.4
$ws.Range("J45").Value = 2487.5
$ws.Range("K45").Value = 1319.4
$ws.Range("L45").Value = 2487.5
$ws.Range("M45").Value = -942.4000000000001
$ws.Range("N45").Value = -3241.5
# Row 63
$ws.Range("H63").Value = 4205.5
$ws.Range("I63").Value = 2712.0386
$ws.Range("J63").Value = 9059.25
$ws.Range("K63").Value = 2712.0386
$ws.Range("L63").Value = 9059.25
$ws.Range("M63").Value = -2026.0386
$ws.Range("N63").Value = -10431.25
# Row 66
$ws.Range("H66").Value = 4205.5
$ws.Range("I66").Value = 2712.0386
$ws.Range("J66").Value = 9059.25
$ws.Range("K66").Value = 13560.193
$ws.Range("L66").Value = 45296.25
$ws.Range("M66").Value = -10128.193
$ws.Range("N66").Value = -52160.25
# Row 96
$ws.Range("H96").Value = 63662.5
$ws.Range("J96").Value = 63662.5
$ws.Range("L96").Value = 63662.5
$ws.Range("N96").Value = -69154.5
# Row 100
$ws.Range("H100").Value = 34999.5
$ws.Range("J100").Value = 34999.5
$ws.Range("L100").Value = 34999.5
$ws.Range("N100").Value = -37163.5
# Row 116
$ws.Range("H116").Value = 1464.0667
$ws.Range("J116").Value = 1505.5
$ws.Range("L116").Value = 1505.5
$ws.Range("N116").Value = -6093.5
# Row 119
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").ClearContents()
$ws.Range("N119").Value = 0

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1464.0667
$ws.Range("J3").Value = 1505.5
$ws.Range("L3").Value = 1505.5
$ws.Range("N3").Value = -1733.5
# Row 134
$ws.Range("H134").Value = 5041.5713
$ws.Range("I134").Value = 5041.5713
$ws.Range("K134").Value = 15124.7139
$ws.Range("M134").Value = -12589.7139

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 19235102
$ws.Range("I31").Value = 3089.5945
$ws.Range("J31").Value = 66674068
$ws.Range("K31").Value = 3089.5945
$ws.Range("L31").Value = 66674068
$ws.Range("M31").Value = -2794.5945
$ws.Range("N31").Value = -66674658
# Row 34
$ws.Range("H34").Value = 19235102
$ws.Range("I34").Value = 3089.5945
$ws.Range("J34").Value = 66674068
$ws.Range("K34").Value = 3089.5945
$ws.Range("L34").Value = 66674068
$ws.Range("M34").Value = -2887.5945
$ws.Range("N34").Value = -66674472
# Row 105
$ws.Range("H105").Value = 15958.875
$ws.Range("I105").Value = 2943.3333
$ws.Range("J105").Value = 55005.5
$ws.Range("K105").Value = 2943.3333
$ws.Range("L105").Value = 55005.5
$ws.Range("M105").Value = -1196.3333
$ws.Range("N105").Value = -58499.5
# Row 122
$ws.Range("H122").Value = 1178.84
$ws.Range("I122").Value = 1225.25
$ws.Range("K122").Value = 3675.75
$ws.Range("M122").Value = -1225.75
# Row 131
$ws.Range("H131").Value = 67438.8
$ws.Range("I131").Value = 35000
$ws.Range("K131").Value = 35000
$ws.Range("M131").Value = -29960

$ws = $wb.Worksheets.Item("CUL")
# Row 130
$ws.Range("H130").Value = 3020
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("M130").ClearContents()
# Row 131
$ws.Range("H131").Value = 36805.25
$ws.Range("J131").Value = 5671.7827
$ws.Range("L131").Value = 17015.3481
$ws.Range("N131").Value = -27095.3481
# Row 132
$ws.Range("H132").Value = 1962665.8
$ws.Range("J132").Value = 5131525.5
$ws.Range("L132").Value = 46183729.5
$ws.Range("N132").Value = -46188789.5
# Row 140
$ws.Range("H140").Value = 2046
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2712.1667
$ws.Range("I80").Value = 2712.1667
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2712.1667
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -1714.1667
# Row 83
$ws.Range("H83").Value = 2712.1667
$ws.Range("I83").Value = 2712.1667
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 13560.8335
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -8568.833500000001
# Row 97
$ws.Range("H97").Value = 1390.0625
$ws.Range("I97").Value = 441
$ws.Range("J97").Value = 2971.8333
$ws.Range("K97").Value = 441
$ws.Range("L97").Value = 2971.8333
$ws.Range("M97").Value = 55
$ws.Range("N97").Value = -3963.8333
# Row 132
$ws.Range("H132").Value = 3284.913
$ws.Range("I132").Value = 3264.6667
$ws.Range("K132").Value = 9794.000100000001
$ws.Range("M132").Value = -7264.000100000001

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2306.3076
$ws.Range("I22").Value = 1565.625
$ws.Range("J22").Value = 2821.5652
$ws.Range("K22").Value = 1565.625
$ws.Range("L22").Value = 2821.5652
$ws.Range("M22").Value = -1270.625
$ws.Range("N22").Value = -3411.5652
# Row 27
$ws.Range("H27").Value = 2306.3076
$ws.Range("I27").Value = 1565.625
$ws.Range("J27").Value = 2821.5652
$ws.Range("K27").Value = 1565.625
$ws.Range("L27").Value = 2821.5652
$ws.Range("M27").Value = -1458.625
$ws.Range("N27").Value = -3035.5652
# Row 55
$ws.Range("H55").Value = 895.64703
$ws.Range("I55").Value = 613.8889
$ws.Range("J55").Value = 1212.625
$ws.Range("K55").Value = 613.8889
$ws.Range("L55").Value = 1212.625
$ws.Range("M55").Value = -440.8889
$ws.Range("N55").Value = -1558.625
# Row 61
$ws.Range("H61").Value = 4658.769
$ws.Range("I61").Value = 3639.8572
$ws.Range("K61").Value = 3639.8572
$ws.Range("M61").Value = -3437.8572
# Row 113
$ws.Range("H113").Value = 4658.769
$ws.Range("I113").Value = 3639.8572
$ws.Range("K113").Value = 3639.8572
$ws.Range("M113").Value = -1469.8572
# Row 122
$ws.Range("H122").Value = 4682.364
$ws.Range("I122").Value = 3813.25
$ws.Range("J122").Value = 7000
$ws.Range("K122").Value = 11439.75
$ws.Range("L122").Value = 21000
$ws.Range("M122").Value = -8989.75
$ws.Range("N122").Value = -25900
# Row 131
$ws.Range("H131").Value = 89239
$ws.Range("J131").Value = 89239
$ws.Range("L131").Value = 89239
$ws.Range("N131").Value = -99319
# Row 132
$ws.Range("H132").Value = 111114980
$ws.Range("I132").Value = 3310.889
$ws.Range("K132").Value = 9932.667000000001
$ws.Range("M132").Value = -7402.667000000001

$ws = $wb.Worksheets.Item("WVR")
# Row 110
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").ClearContents()
$ws.Range("N110").Value = 0
# Row 122
$ws.Range("H122").Value = 71502010
$ws.Range("I122").Value = 83418264
$ws.Range("K122").Value = 250254792
$ws.Range("M122").Value = -250252342
# Row 124
$ws.Range("H124").Value = 50000
$ws.Range("J124").Value = 50000
$ws.Range("L124").Value = 50000
$ws.Range("N124").Value = -59820
# Row 132
$ws.Range("H132").Value = 5774.5713
$ws.Range("I132").Value = 5617.5806
$ws.Range("J132").Value = 6991.25
$ws.Range("K132").Value = 16852.7418
$ws.Range("L132").Value = 20973.75
$ws.Range("M132").Value = -14322.7418
$ws.Range("N132").Value = -26033.75

Write-Output "done"